# Great Circle Calculations Added
# Adds Lat_Rad / Long_Rad / Delta / Distance columns next to the existing
# Top-Left / Bottom-Right lat-long table and computes the great-circle
# distance between the two points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (D1:G1) ---------------------------------------------------
$ws.Range("D1").Value = "Lat_Rad"
$ws.Range("E1").Value = "Long_Rad"
$ws.Range("F1").Value = "Delta"
$ws.Range("G1").Value = "Distance"
$ws.Range("G1").Font.Bold = $true

# --- Helper columns: convert latitude/longitude to radians ---------------
$ws.Range("D2").Formula = "=RADIANS(B2)"
$ws.Range("D3").Formula = "=RADIANS(B3)"

# Long_Rad is entered as a (legacy CSE-style) array formula spilling over
# E2:E3, matching RADIANS(C2:C3).
$ws.Range("E2:E3").FormulaArray = "=RADIANS(C2:C3)"

# Difference in longitude (radians)
$ws.Range("F2").Formula = "=E3-E2"

# Great-circle distance (km), using the radius of the Earth in B5
$ws.Range("G2").Formula = "=ACOS(SIN(D2)*SIN(D3)+COS(D2)*COS(D3)*COS(F2))*B5"
$ws.Range("G2").Font.Bold = $true

# --- Cosmetics -------------------------------------------------------------
# Widen the Long_Rad column enough to comfortably fit its values.
$ws.Columns("E").ColumnWidth = 8.6

# Portrait page orientation for printing.
$ws.PageSetup.Orientation = 1

# Leave the selection on the newly-computed distance cell.
$ws.Range("G2").Select() | Out-Null
